# "Update (Removed Auto Arima)"
# Updates forecast numbers on the "Forecast Comparison" sheet (Prophet / Amazon
# Mean / P70 / P90 forecast columns) and the rolled-up totals on the "Summary"
# sheet, reflecting removal of the Auto-ARIMA model from the forecast blend.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: Prophet Forecast (C), Amazon Mean Forecast (D),
#     Amazon P70 Forecast (E), Amazon P80 Forecast (F), Amazon P90 Forecast (G)

$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 2

$ws1.Range("D3").Value = 2
$ws1.Range("E3").Value = 2
$ws1.Range("G3").Value = 4

$ws1.Range("C4").Value = 2
$ws1.Range("D4").Value = 2

$ws1.Range("C5").Value = 2
$ws1.Range("D5").Value = 2

$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 2
$ws1.Range("E6").Value = 2
$ws1.Range("G6").Value = 4

$ws1.Range("C7").Value = 3
$ws1.Range("D7").Value = 2
$ws1.Range("E7").Value = 2
$ws1.Range("F7").Value = 3
$ws1.Range("G7").Value = 4

$ws1.Range("C8").Value = 3
$ws1.Range("D8").Value = 2
$ws1.Range("G8").Value = 5

$ws1.Range("C9").Value = 3
$ws1.Range("D9").Value = 2

$ws1.Range("C10").Value = 2
$ws1.Range("D10").Value = 2
$ws1.Range("E10").Value = 2
$ws1.Range("F10").Value = 3

$ws1.Range("C11").Value = 2
$ws1.Range("D11").Value = 2
$ws1.Range("F11").Value = 3
$ws1.Range("G11").Value = 5

$ws1.Range("C12").Value = 1
$ws1.Range("D12").Value = 2

$ws1.Range("D13").Value = 2

$ws1.Range("D14").Value = 2

$ws1.Range("C15").Value = 1
$ws1.Range("D15").Value = 2
$ws1.Range("G15").Value = 6

$ws1.Range("C16").Value = 2
$ws1.Range("D16").Value = 2
$ws1.Range("G16").Value = 6

$ws1.Range("C17").Value = 1
$ws1.Range("D17").Value = 2
$ws1.Range("G17").Value = 6

# --- Summary sheet totals. These are stored as plain text (not numbers /
#     dates), so force text formatting before assigning so the numeric- and
#     date-looking strings aren't auto-converted to a number / date serial.

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "30"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "19"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "7"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "3"

$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2024-12-08"
